$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 2469585.2
$ws.Range("I33").Value = 4629888.5
$ws.Range("J33").Value = 667
$ws.Range("K33").Value = 4629888.5
$ws.Range("L33").Value = 667
$ws.Range("M33").Value = -4629659.5
$ws.Range("N33").Value = -1125

# Row 70
$ws.Range("H70").Value = 7379.8
$ws.Range("I70").Value = 1800
$ws.Range("J70").Value = 7999.778
$ws.Range("K70").Value = 5400
$ws.Range("L70").Value = 23999.334
$ws.Range("M70").Value = -5130
$ws.Range("N70").Value = -24539.334

# Row 73
$ws.Range("H73").Value = 7379.8
$ws.Range("I73").Value = 1800
$ws.Range("J73").Value = 7999.778
$ws.Range("K73").Value = 5400
$ws.Range("L73").Value = 23999.334
$ws.Range("M73").Value = -4464
$ws.Range("N73").Value = -25871.334

# Row 74
$ws.Range("H74").Value = 7433.2334
$ws.Range("J74").Value = 7883.9165
$ws.Range("L74").Value = 7883.9165
$ws.Range("N74").Value = -9755.916499999999

# Row 77
$ws.Range("H77").Value = 7433.2334
$ws.Range("J77").Value = 7883.9165
$ws.Range("L77").Value = 39419.5825
$ws.Range("N77").Value = -48779.5825

# Row 93
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -39992

# Row 96
$ws.Range("H96").Value = 228.92857
$ws.Range("I96").Value = 231.54546
$ws.Range("K96").Value = 694.6363799999999
$ws.Range("M96").Value = 678.3636200000001

# Row 101
$ws.Range("H101").Value = 38462016
$ws.Range("I101").Value = 50000440
$ws.Range("J101").Value = 595
$ws.Range("K101").Value = 150001320
$ws.Range("L101").Value = 1785
$ws.Range("M101").Value = -149999698
$ws.Range("N101").Value = -5029

# Row 131
$ws.Range("H131").Value = 3569.3076
$ws.Range("I131").Value = 1177.7059
$ws.Range("K131").Value = 3533.1177
$ws.Range("M131").Value = 1506.8823

# Row 135
$ws.Range("H135").Value = 1056.7693
$ws.Range("J135").Value = 1829.4286
$ws.Range("L135").Value = 16464.8574
$ws.Range("N135").Value = -21534.8574

# Row 137
$ws.Range("H137").Value = 76135.625
$ws.Range("I137").Value = 105969.12
$ws.Range("K137").Value = 317907.36
$ws.Range("M137").Value = -315357.36

# Row 138
$ws.Range("H138").Value = 3674.9473
$ws.Range("J138").Value = 3753.52
$ws.Range("L138").Value = 11260.56
$ws.Range("N138").Value = -21540.56

# Row 141
$ws.Range("H141").Value = 14271.615
$ws.Range("I141").Value = 7506.773
$ws.Range("K141").Value = 22520.319
$ws.Range("M141").Value = -17340.319

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1886030.9
$ws.Range("I2").Value = 2356960.5
$ws.Range("K2").Value = 2356960.5
$ws.Range("M2").Value = -2356847.5

# Row 32
$ws.Range("H32").Value = 9490.661
$ws.Range("I32").Value = 5536.146
$ws.Range("J32").Value = 23049
$ws.Range("K32").Value = 5536.146
$ws.Range("L32").Value = 23049
$ws.Range("M32").Value = -5249.146
$ws.Range("N32").Value = -23623

# Row 63
$ws.Range("H63").Value = 6991.6665
$ws.Range("I63").Value = 2831.6667
$ws.Range("J63").Value = 8378.333000000001
$ws.Range("K63").Value = 2831.6667
$ws.Range("L63").Value = 8378.333000000001
$ws.Range("M63").Value = -2145.6667
$ws.Range("N63").Value = -9750.333000000001

# Row 66
$ws.Range("H66").Value = 6991.6665
$ws.Range("I66").Value = 2831.6667
$ws.Range("J66").Value = 8378.333000000001
$ws.Range("K66").Value = 14158.3335
$ws.Range("L66").Value = 41891.665
$ws.Range("M66").Value = -10726.3335
$ws.Range("N66").Value = -48755.665

# Row 103
$ws.Range("H103").Value = 60177.5
$ws.Range("J103").Value = 60177.5
$ws.Range("L103").Value = 60177.5
$ws.Range("N103").Value = -62521.5

# Row 116
$ws.Range("H116").Value = 1886030.9
$ws.Range("I116").Value = 2356960.5
$ws.Range("K116").Value = 2356960.5
$ws.Range("M116").Value = -2354666.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1886030.9
$ws.Range("I3").Value = 2356960.5
$ws.Range("K3").Value = 2356960.5
$ws.Range("M3").Value = -2356846.5

# Row 86
$ws.Range("H86").Value = 5558213
$ws.Range("I86").Value = 7145684.5
$ws.Range("J86").Value = 2062.5
$ws.Range("K86").Value = 7145684.5
$ws.Range("L86").Value = 2062.5
$ws.Range("M86").Value = -7144561.5
$ws.Range("N86").Value = -4308.5

# Row 89
$ws.Range("H89").Value = 5558213
$ws.Range("I89").Value = 7145684.5
$ws.Range("J89").Value = 2062.5
$ws.Range("K89").Value = 35728422.5
$ws.Range("L89").Value = 10312.5
$ws.Range("M89").Value = -35722806.5
$ws.Range("N89").Value = -21544.5

# Row 103
$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344

# Row 134
$ws.Range("H134").Value = 5114.3
$ws.Range("I134").Value = 1989.2727
$ws.Range("K134").Value = 5967.8181
$ws.Range("M134").Value = -3432.8181

# Row 139
$ws.Range("H139").Value = 134744.38
$ws.Range("J139").Value = 147992.67
$ws.Range("L139").Value = 147992.67
$ws.Range("N139").Value = -158272.67

# Row 140
$ws.Range("H140").Value = 88779
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 88779
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 88779
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -99139

$ws = $wb.Worksheets.Item("CRP")
# Row 43
$ws.Range("H43").Value = 29614.25
$ws.Range("J43").Value = 29614.25
$ws.Range("L43").Value = 29614.25
$ws.Range("N43").Value = -29982.25

# Row 101
$ws.Range("H101").Value = 29614.25
$ws.Range("J101").Value = 29614.25
$ws.Range("L101").Value = 29614.25
$ws.Range("N101").Value = -36104.25

# Row 134
$ws.Range("H134").Value = 4089.9375
$ws.Range("I134").Value = 2501.7144
$ws.Range("J134").Value = 5325.222
$ws.Range("K134").Value = 7505.1432
$ws.Range("L134").Value = 15975.666
$ws.Range("M134").Value = -4970.1432
$ws.Range("N134").Value = -21045.666

# Row 141
$ws.Range("H141").Value = 137934.7
$ws.Range("J141").Value = 137934.7
$ws.Range("L141").Value = 137934.7
$ws.Range("N141").Value = -148294.7

$ws = $wb.Worksheets.Item("CUL")
# Row 141
$ws.Range("I141").Value = 2859.6
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8578.799999999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3398.799999999999
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 76114.664
$ws.Range("J95").Value = 76114.664
$ws.Range("L95").Value = 76114.664
$ws.Range("N95").Value = -81606.664

# Row 102
$ws.Range("H102").Value = 8997968
$ws.Range("I102").Value = 22223400
$ws.Range("K102").Value = 22223400
$ws.Range("M102").Value = -22221778

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 38672.582
$ws.Range("I22").Value = 60145.2
$ws.Range("J22").Value = 2884.889
$ws.Range("K22").Value = 60145.2
$ws.Range("L22").Value = 2884.889
$ws.Range("M22").Value = -59850.2
$ws.Range("N22").Value = -3474.889

# Row 27
$ws.Range("H27").Value = 38672.582
$ws.Range("I27").Value = 60145.2
$ws.Range("J27").Value = 2884.889
$ws.Range("K27").Value = 60145.2
$ws.Range("L27").Value = 2884.889
$ws.Range("M27").Value = -60038.2
$ws.Range("N27").Value = -3098.889

# Row 46
$ws.Range("H46").Value = 4356724.5
$ws.Range("J46").Value = 9887
$ws.Range("L46").Value = 9887
$ws.Range("N46").Value = -10263

# Row 55
$ws.Range("H55").Value = 1931.625
$ws.Range("J55").Value = 302.75
$ws.Range("L55").Value = 302.75
$ws.Range("N55").Value = -648.75

# Row 64
$ws.Range("H64").Value = 500015070
$ws.Range("J64").Value = 500015070
$ws.Range("L64").Value = 500015070
$ws.Range("N64").Value = -500015520

# Row 67
$ws.Range("H67").Value = 500015070
$ws.Range("J67").Value = 500015070
$ws.Range("L67").Value = 500015070
$ws.Range("N67").Value = -500016630

# Row 93
$ws.Range("H93").Value = 47648236
$ws.Range("I93").Value = 111111550
$ws.Range("J93").Value = 50750.75
$ws.Range("K93").Value = 111111550
$ws.Range("L93").Value = 50750.75
$ws.Range("M93").Value = -111110302
$ws.Range("N93").Value = -53246.75

# Row 103
$ws.Range("H103").Value = 28333
$ws.Range("J103").Value = 28333
$ws.Range("L103").Value = 28333
$ws.Range("N103").Value = -30677

# Row 122
$ws.Range("H122").Value = 5764.0435
$ws.Range("I122").Value = 3421.111
$ws.Range("J122").Value = 7270.2144
$ws.Range("K122").Value = 10263.333
$ws.Range("L122").Value = 21810.6432
$ws.Range("M122").Value = -7813.332999999999
$ws.Range("N122").Value = -26710.6432

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 55549
$ws.Range("J63").Value = 55549
$ws.Range("L63").Value = 55549
$ws.Range("N63").Value = -56797

# Row 66
$ws.Range("H66").Value = 55549
$ws.Range("J66").Value = 55549
$ws.Range("L66").Value = 166647
$ws.Range("N66").Value = -172887

# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
